$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update attendance counters from 0 to 1 for the relevant cells
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

$ws.Range("H7").Value = 1

$ws.Range("H8").Value = 1

$ws.Range("H9").Value = 1

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

$ws.Range("H11").Value = 1

$ws.Range("H12").Value = 1

$ws.Range("H13").Value = 1

$ws.Range("H14").Value = 1

$ws.Range("H15").Value = 1

$ws.Range("H16").Value = 1

$ws.Range("H17").Value = 1

$ws.Range("H18").Value = 1
